# CargadorManzanas y visualizacion manzanas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (columns swapped: B becomes Coordenadas, C becomes Nro manzanas vecinas) ---
$ws.Range("B1").Value = "Coordenadas"
$ws.Range("C1").Value = "Nro manzanas vecinas"

# --- Column widths ---
# NOTE: the engine's ColumnWidth setter only has ~1/6-character (pixel-level)
# granularity, so it cannot reproduce fractional widths such as
# 40.42578125 / 25.140625 exactly; these inputs land mid-band on the closest
# achievable stored width (40.5 / 25.16666...) for robustness.
$ws.Columns("B").ColumnWidth = 39.65
$ws.Columns("C").ColumnWidth = 17.15
$ws.Columns("D").ColumnWidth = 24.3

# --- Coordinates (column B) ---
$ws.Range("B2").Value = "-34.51909573235544, -58.71988115476912"
$ws.Range("B3").Value = "-34.52637948225431, -58.71009645709745"
$ws.Range("B4").Value = "-34.528005863683184, -58.731039143342066"
$ws.Range("B5").Value = "-34.54045016573941, -58.72511682633027"
$ws.Range("B6").Value = "-34.54589396298188, -58.715589620702595"
$ws.Range("B7").Value = "-34.53599588498832, -58.69670687080991"
$ws.Range("B8").Value = "-34.553275292650206, -58.719442876513725"
$ws.Range("B9").Value = "-34.56213183370967, -58.70717317543131"
$ws.Range("B10").Value = "-34.5586715072938, -58.691521521470904"

# --- Neighbouring-block counts (column C) ---
$ws.Range("C2").Value = 2.1
$ws.Range("C3").Value = 0.5
$ws.Range("C4").Value = 0.3
$ws.Range("C5").Value = "2,4,6"
$ws.Range("C6").Value = "3,5,7"
$ws.Range("C7").Value = "1,4,8"
$ws.Range("C8").Value = 3.7
$ws.Range("C9").Value = "4,6,8"
$ws.Range("C10").Value = 5.7

# --- Apply Text format AFTER populating so the numeric cells stay numeric ---
# (matches Excel behaviour: formatting a cell as Text does not retype already
# entered numbers, it only affects future keystrokes/parsing)
$ws.Range("B2:C10").NumberFormat = "@"

# --- New row 10 block number ---
$ws.Range("A10").Value = 8

# --- Selection / view state ---
[void]$ws.Range("B2:C10").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
